$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Suspect sheet: remove the old fk_CrimeID / fk_StateID columns (M, N).
#    Those one-to-one relations are being replaced by the new many-to-many
#    WantedFor / WantedIn sheets below.
# ---------------------------------------------------------------------------
$suspect = $wb.Worksheets.Item("Suspect")
$suspect.Range("M1:N4").ClearContents()

# ---------------------------------------------------------------------------
# 2. State sheet: add a Population column (E).
# ---------------------------------------------------------------------------
$state = $wb.Worksheets.Item("State")
$state.Range("E1").Value = "Population"

$populations = @(4863300,741894,6931071,2988248,39250017,5540545,3576452,952065,20612439,10310371,1428557,1683140,12801539,6633053,3134693,2907289,4436974,4681666,1331479,6016447,6811779,9928300,5519952,2988726,6093000,1042520,1907116,2940058,1334795,8944469,2081015,19745289,10146788,757952,11614373,3923561,4093465,12784227,1056426,4961119,865454,6651194,27862596,3051217,624594,8411808,7288000,681170,1831102,5778708,585501)

for ($i = 0; $i -lt $populations.Length; $i++) {
    $row = $i + 2
    $state.Cells.Item($row, 5).Value = $populations[$i]
}

# ---------------------------------------------------------------------------
# 3. New sheet "WantedIn": fk_SuspectID -> fk_StateID (many-to-many).
#    Inserted right after "Crime" (before "CrimeByState").
# ---------------------------------------------------------------------------
$crime = $wb.Worksheets.Item("Crime")
$wantedIn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $crime)
$wantedIn.Name = "WantedIn"

$wantedIn.Range("A1").Value = "fk_SuspectID"
$wantedIn.Range("B1").Value = "fk_StateID"

$wantedInData = @(
    @(1, 37),
    @(1, 47),
    @(1, 5),
    @(1, 6),
    @(1, 51),
    @(2, 35),
    @(3, 20)
)

for ($i = 0; $i -lt $wantedInData.Length; $i++) {
    $row = $i + 2
    $wantedIn.Cells.Item($row, 1).Value = $wantedInData[$i][0]
    $wantedIn.Cells.Item($row, 2).Value = $wantedInData[$i][1]
}

# ---------------------------------------------------------------------------
# 4. New sheet "WantedFor": fk_SuspectID -> fk_CrimeID (many-to-many).
#    Inserted right after "WantedIn" (before "CrimeByState").
# ---------------------------------------------------------------------------
$wantedFor = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wantedIn)
$wantedFor.Name = "WantedFor"

$wantedFor.Range("A1").Value = "fk_SuspectID"
$wantedFor.Range("B1").Value = "fk_CrimeID"

$wantedForData = @(
    @(1, 1),
    @(2, 2),
    @(3, 3)
)

for ($i = 0; $i -lt $wantedForData.Length; $i++) {
    $row = $i + 2
    $wantedFor.Cells.Item($row, 1).Value = $wantedForData[$i][0]
    $wantedFor.Cells.Item($row, 2).Value = $wantedForData[$i][1]
}
